$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "darsh2@gmail.com"
$ws.Range("B2").Value = "sanj2@gmail.com"
$ws.Range("B3").Value = "harshi2@gmail.com"
